$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the three anchor paragraphs *before* making any changes, so that
# paragraph indices found here are still correct once we start inserting
# new paragraphs (we apply the edits from the bottom of the document
# upwards below, so earlier indices never shift under us).
# ---------------------------------------------------------------------------
function Find-ParagraphContaining($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.Contains($text)) {
            return $i
        }
    }
    return -1
}

$variusPara = Find-ParagraphContaining "[Varius][Artorius][Musai][Anaxu][Valia] - [Eresse][Tyrios][Faros]"
$abyzuPara  = Find-ParagraphContaining "[Abyzu][Daria][Varius][Eressea][Galilon][Dundain][Eldar][Laniakaea][Vathrin]"
$tyrPara    = Find-ParagraphContaining "[Tyr] |::|"

if ($variusPara -eq -1 -or $abyzuPara -eq -1 -or $tyrPara -eq -1) {
    throw "Could not locate one or more anchor paragraphs (Varius=$variusPara, Abyzu=$abyzuPara, Tyr=$tyrPara)"
}

# ---------------------------------------------------------------------------
# Edit 3 (latest in the doc -> applied first): the paragraph "[Tyr] |::|"
# gains extra bracketed content, and a new paragraph "[Dalphine][Sirion]"
# is inserted right after it.
# ---------------------------------------------------------------------------
$d.Paragraphs($tyrPara).Range.InsertAfter(" [Megalon][Antioch][Antorus][Altair] | [Asphodel][Artorius][Arcadius][Anorius]")

$d.Paragraphs($tyrPara).Range.InsertParagraphAfter()
$newRange = $d.Paragraphs($tyrPara + 1).Range
$d.Range($newRange.Start, $newRange.End - 1).Text = "[Dalphine][Sirion]"

# ---------------------------------------------------------------------------
# Edit 2: a new paragraph "[Megalon][Dalphine][Regulus][Tyrion][Eldar][Eresse]"
# is inserted right after the "[Abyzu][Daria]...[Vathrin]" paragraph.
# ---------------------------------------------------------------------------
$d.Paragraphs($abyzuPara).Range.InsertParagraphAfter()
$newRange2 = $d.Paragraphs($abyzuPara + 1).Range
$d.Range($newRange2.Start, $newRange2.End - 1).Text = "[Megalon][Dalphine][Regulus][Tyrion][Eldar][Eresse]"

# ---------------------------------------------------------------------------
# Edit 1 (earliest in the doc -> applied last): the first of the run of
# empty paragraphs right after "[Varius][Artorius][Musai][Anaxu][Valia] -
# [Eresse][Tyrios][Faros]" gets text "[Antorus][Megalon][Artorius][Musai]",
# and a new empty paragraph is inserted right after it.
# ---------------------------------------------------------------------------
$targetIndex = $variusPara + 1
$d.Paragraphs($targetIndex).Range.InsertBefore("[Antorus][Megalon][Artorius][Musai]")
$d.Paragraphs($targetIndex).Range.InsertParagraphAfter()
